# Updates the "北京-漫展信息" workbook to the refreshed scrape output
# (gh-pages output generated at dd351a1).
#
# Changes applied:
#  - "想去人数" (F column) refreshed to newer counts on several existing rows
#    across all four sheets (展览 / 演出 / 本地生活 / 全部类型).
#  - On 展览 (sheet1) and 全部类型 (sheet4) the 2024.03.16/2024.03.23 block
#    was re-scraped: the bilibili event id=78902 changed its name/date from
#    "北京·yiyou 双马 运动番动漫展" (2024.03.23) to
#    "北京· YiYou 运动番only 2.0" (2024.03.16), which re-sorts it to the top
#    of that date block, pushing the other two rows down by one position.

$wb = $excel.ActiveWorkbook

function Set-TextCell($ws, $cellRef, $text) {
    # Force text storage even when the string looks numeric (e.g. "60"),
    # matching the source data's column G ("最低票价") which is text-typed.
    $ws.Range($cellRef).NumberFormat = "@"
    $ws.Range($cellRef).Value = $text
}

# ---------------------------------------------------------------------
# Sheet "展览"
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("展览")

$ws1.Range("F5").Value = 1978
$ws1.Range("F7").Value = 435
$ws1.Range("F9").Value = 213
$ws1.Range("F10").Value = 6942
$ws1.Range("F12").Value = 542
$ws1.Range("F13").Value = 118
$ws1.Range("F15").Value = 2397
$ws1.Range("F16").Value = 1747
$ws1.Range("F18").Value = 48
$ws1.Range("F19").Value = 100
$ws1.Range("F21").Value = 114

# Row 23: now the re-scraped "YiYou 运动番only 2.0" event (was at row 25,
# id=78902, date moved from 2024.03.23 to 2024.03.16).
Set-TextCell $ws1 "B23" "2024.03.16"
$ws1.Range("C23").Value = "北京· YiYou 运动番only 2.0"
$ws1.Range("D23").Value = "京开高速入口与京开高速交叉口西180米 北京双马文体创业园"
$ws1.Range("E23").Value = "2024.03.16 09:30-03.17 18:00"
$ws1.Range("F23").Value = 168
Set-TextCell $ws1 "G23" "60"
$ws1.Range("H23").Value = $False
$ws1.Range("I23").Value = "https://show.bilibili.com/platform/detail.html?id=78902&msource=Msearch_colligation"

# Row 24: shifted down from the old row 23 (thebONE GOJO), F refreshed.
Set-TextCell $ws1 "B24" "2024.03.16"
$ws1.Range("C24").Value = "北京·thebONE✖️GOJO超次元嘉年华02"
$ws1.Range("D24").Value = "小关路39号 北投购物公园"
$ws1.Range("E24").Value = "2024.03.16 10:00-03.17 17:00"
$ws1.Range("F24").Value = 81
Set-TextCell $ws1 "G24" "70"
$ws1.Range("H24").Value = $False
$ws1.Range("I24").Value = "https://show.bilibili.com/platform/detail.html?id=78896&msource=Msearch_colligation"

# Row 25: shifted down from the old row 24 (thebONE x Ilike), F refreshed.
Set-TextCell $ws1 "B25" "2024.03.23"
$ws1.Range("C25").Value = "北京·thebONE×Ilike动漫游戏嘉年华S4"
$ws1.Range("D25").Value = "小关路39号 北投购物公园"
$ws1.Range("E25").Value = "2024.03.23 10:00-03.24 17:00"
$ws1.Range("F25").Value = 975
Set-TextCell $ws1 "G25" "58"
$ws1.Range("H25").Value = $True
$ws1.Range("I25").Value = "https://show.bilibili.com/platform/detail.html?id=79601&msource=Msearch_colligation"

$ws1.Range("F26").Value = 155
$ws1.Range("F27").Value = 4108

# ---------------------------------------------------------------------
# Sheet "演出"
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 71

# ---------------------------------------------------------------------
# Sheet "本地生活"
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F3").Value = 688

# ---------------------------------------------------------------------
# Sheet "全部类型" (aggregated view, same edits as 展览 + the other two
# sheets, offset by +5 rows in the 2024.03.16/2024.03.23 block).
# ---------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("全部类型")

$ws4.Range("F4").Value = 688
$ws4.Range("F6").Value = 71
$ws4.Range("F8").Value = 1978
$ws4.Range("F12").Value = 435
$ws4.Range("F14").Value = 213
$ws4.Range("F15").Value = 6942
$ws4.Range("F17").Value = 542
$ws4.Range("F18").Value = 118
$ws4.Range("F20").Value = 2397
$ws4.Range("F21").Value = 1747
$ws4.Range("F23").Value = 48
$ws4.Range("F24").Value = 100
$ws4.Range("F26").Value = 114

# Row 28: re-scraped "YiYou 运动番only 2.0" event (was at row 30).
Set-TextCell $ws4 "B28" "2024.03.16"
$ws4.Range("C28").Value = "北京· YiYou 运动番only 2.0"
$ws4.Range("D28").Value = "京开高速入口与京开高速交叉口西180米 北京双马文体创业园"
$ws4.Range("E28").Value = "2024.03.16 09:30-03.17 18:00"
$ws4.Range("F28").Value = 168
Set-TextCell $ws4 "G28" "60"
$ws4.Range("H28").Value = $False
$ws4.Range("I28").Value = "https://show.bilibili.com/platform/detail.html?id=78902&msource=Msearch_colligation"

# Row 29: shifted down from the old row 28 (thebONE GOJO), F refreshed.
Set-TextCell $ws4 "B29" "2024.03.16"
$ws4.Range("C29").Value = "北京·thebONE✖️GOJO超次元嘉年华02"
$ws4.Range("D29").Value = "小关路39号 北投购物公园"
$ws4.Range("E29").Value = "2024.03.16 10:00-03.17 17:00"
$ws4.Range("F29").Value = 81
Set-TextCell $ws4 "G29" "70"
$ws4.Range("H29").Value = $False
$ws4.Range("I29").Value = "https://show.bilibili.com/platform/detail.html?id=78896&msource=Msearch_colligation"

# Row 30: shifted down from the old row 29 (thebONE x Ilike), F refreshed.
Set-TextCell $ws4 "B30" "2024.03.23"
$ws4.Range("C30").Value = "北京·thebONE×Ilike动漫游戏嘉年华S4"
$ws4.Range("D30").Value = "小关路39号 北投购物公园"
$ws4.Range("E30").Value = "2024.03.23 10:00-03.24 17:00"
$ws4.Range("F30").Value = 975
Set-TextCell $ws4 "G30" "58"
$ws4.Range("H30").Value = $True
$ws4.Range("I30").Value = "https://show.bilibili.com/platform/detail.html?id=79601&msource=Msearch_colligation"

$ws4.Range("F31").Value = 155
$ws4.Range("F32").Value = 4108
